$d = $word.ActiveDocument

# Locate the "Requisitos" Heading2 paragraph that starts the block to remove
# (that heading plus the "LOQ4205 ... (Requisito fraco)" bullet that follows
# it, running to the end of the document body) and delete everything from
# the start of that paragraph through the end of the document's content.
$requisitosStart = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Requisitos") {
        $requisitosStart = $p.Range.Start
        break
    }
}

if ($requisitosStart -ge 0) {
    $r = $d.Range($requisitosStart, $d.Content.End)
    $r.Delete()
}
